# Weekly data refresh for "Fruta, Vega Modelo de Temuco - Pera asiática".
# The rows keep the same header/shape; only the per-record data (date,
# quality, volume, prices, unit of sale, origin, $/Kg, Kg/unit) is
# refreshed with a newer weekly snapshot. Some rows effectively swap
# their data with other rows (the same record landing on a different
# row) while rows 13, 15 and 17 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data that moves between rows.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Mapping: destination row -> source row (source row's data, as it was
# BEFORE this edit, ends up in the destination row).
$rowMap = @{
    2  = 20
    3  = 12
    4  = 9
    5  = 10
    6  = 8
    7  = 5
    8  = 2
    9  = 3
    10 = 7
    11 = 6
    12 = 18
    14 = 11
    16 = 14
    18 = 19
    19 = 16
    20 = 4
}

# Snapshot every source row's current values before writing anything,
# since several rows both give and receive data (cycles in the
# permutation), so in-place writes would clobber values still needed.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($col in $cols) {
            # Value2 round-trips numbers/strings/dates cleanly; plain
            # Value returns a COM variant wrapper here that doesn't.
            $rowData[$col] = $ws.Range("$col$srcRow").Value2
        }
        $snapshot[$srcRow] = $rowData
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}
